# Update crypto price/volume data per the Tue Mar  5 08:36:21 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.480.49'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").Value = '3.683.60'
$ws.Range("E3").Value = '  +4.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '417.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("D7").Value = '3.679.02'
$ws.Range("E7").Value = '  +4.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.638'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.756'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.179'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000384'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +41.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.17%  '
$ws.Range("D15").Value = '4.271.13'
$ws.Range("E15").Value = '  +5.06%  '
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("D17").Value = '3.801.05'
$ws.Range("E17").Value = '  +7.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '66.631.49'
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '437.33'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '16.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +20.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("E25").Value = '  -5.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.06%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.124'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.29%  '
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.19'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0491'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +33.99%  '
$ws.Range("D40").Value = '0.0₃0722'
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("E41").Value = '  +2.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '28.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +29.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.37%  '
$ws.Range("E49").Value = '  -7.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.303'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.54%  '
$ws.Range("E51").Value = '  +10.45%  '
